$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3014.2104
$ws.Range("J40").Value = 2176
$ws.Range("L40").Value = 2176
$ws.Range("N40").Value = -2526
$ws.Range("H58").Value = 935.3125
$ws.Range("I58").Value = 458.84616
$ws.Range("K58").Value = 1376.53848
$ws.Range("M58").Value = -1226.53848
$ws.Range("H82").Value = 507.25
$ws.Range("I82").Value = 507.25
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1521.75
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1115.75
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 507.25
$ws.Range("I85").Value = 507.25
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1521.75
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -117.75
$ws.Range("N85").ClearContents()
$ws.Range("H92").Value = 566.9286
$ws.Range("I92").Value = 451.36365
$ws.Range("J92").Value = 990.6667
$ws.Range("K92").Value = 451.36365
$ws.Range("L92").Value = 990.6667
$ws.Range("M92").Value = 796.63635
$ws.Range("N92").Value = -3486.6667
$ws.Range("H98").Value = 749.8261
$ws.Range("I98").Value = 770.4761999999999
$ws.Range("K98").Value = 770.4761999999999
$ws.Range("M98").Value = 727.5238000000001
$ws.Range("H118").Value = 5243.5
$ws.Range("I118").Value = 783.6923
$ws.Range("J118").Value = 8653.941000000001
$ws.Range("K118").Value = 2351.0769
$ws.Range("L118").Value = 25961.823
$ws.Range("M118").Value = -694.0769
$ws.Range("N118").Value = -29275.823
$ws.Range("H122").Value = 749.8261
$ws.Range("I122").Value = 770.4761999999999
$ws.Range("K122").Value = 2311.4286
$ws.Range("M122").Value = 138.5714000000003
$ws.Range("H132").Value = 2909.1838
$ws.Range("I132").Value = 2672.9595
$ws.Range("J132").Value = 4253.846
$ws.Range("K132").Value = 8018.8785
$ws.Range("L132").Value = 12761.538
$ws.Range("M132").Value = -5488.8785
$ws.Range("N132").Value = -17821.538
$ws.Range("H137").Value = 38190.85
$ws.Range("I137").Value = 902.8333
$ws.Range("J137").Value = 112766.89
$ws.Range("K137").Value = 2708.4999
$ws.Range("L137").Value = 338300.67
$ws.Range("M137").Value = -158.4998999999998
$ws.Range("N137").Value = -343400.67
$ws.Range("H138").Value = 5365.93
$ws.Range("I138").Value = 1852.9333
$ws.Range("J138").Value = 5985.8706
$ws.Range("K138").Value = 5558.7999
$ws.Range("L138").Value = 17957.6118
$ws.Range("M138").Value = -418.7999
$ws.Range("N138").Value = -28237.6118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34862.414
$ws.Range("I32").Value = 25652.549
$ws.Range("J32").Value = 44647.895
$ws.Range("K32").Value = 25652.549
$ws.Range("L32").Value = 44647.895
$ws.Range("M32").Value = -25365.549
$ws.Range("N32").Value = -45221.895
$ws.Range("H74").Value = 373731.28
$ws.Range("I74").Value = 3238.7917
$ws.Range("J74").Value = 3337671.2
$ws.Range("K74").Value = 3238.7917
$ws.Range("L74").Value = 3337671.2
$ws.Range("M74").Value = -2364.7917
$ws.Range("N74").Value = -3339419.2
$ws.Range("H77").Value = 373731.28
$ws.Range("I77").Value = 3238.7917
$ws.Range("J77").Value = 3337671.2
$ws.Range("K77").Value = 16193.9585
$ws.Range("L77").Value = 16688356
$ws.Range("M77").Value = -11825.9585
$ws.Range("N77").Value = -16697092
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 1584
$ws.Range("I54").Value = 1584
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1584
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1100
$ws.Range("N54").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 2519.6667
$ws.Range("I41").Value = 2519.6667
$ws.Range("K41").Value = 2519.6667
$ws.Range("M41").Value = -2091.6667
$ws.Range("H50").Value = 18533.334
$ws.Range("J50").Value = 18533.334
$ws.Range("L50").Value = 18533.334
$ws.Range("N50").Value = -19783.334
$ws.Range("H51").Value = 12966.333
$ws.Range("J51").Value = 12966.333
$ws.Range("L51").Value = 12966.333
$ws.Range("N51").Value = -14438.333
$ws.Range("H58").Value = 3093.7778
$ws.Range("I58").Value = 628.7059
$ws.Range("K58").Value = 628.7059
$ws.Range("M58").Value = -425.7059
$ws.Range("H59").Value = 19679.666
$ws.Range("J59").Value = 19679.666
$ws.Range("L59").Value = 19679.666
$ws.Range("N59").Value = -21969.666
$ws.Range("H60").Value = 12950
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 12950
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 12950
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -13972
$ws.Range("H61").Value = 12966.333
$ws.Range("J61").Value = 12966.333
$ws.Range("L61").Value = 12966.333
$ws.Range("N61").Value = -13662.333
$ws.Range("H132").Value = 1860.6818
$ws.Range("I132").Value = 1163.1666
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 3489.4998
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -959.4998000000001
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 3093.7778
$ws.Range("I136").Value = 628.7059
$ws.Range("K136").Value = 1886.1177
$ws.Range("M136").Value = 663.8822999999998
$ws.Range("H141").Value = 51340.965
$ws.Range("I141").Value = 30177.8
$ws.Range("J141").Value = 55941.652
$ws.Range("K141").Value = 30177.8
$ws.Range("L141").Value = 55941.652
$ws.Range("M141").Value = -24997.8
$ws.Range("N141").Value = -66301.652
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 684.7049
$ws.Range("I5").Value = 516.44684
$ws.Range("J5").Value = 1249.5714
$ws.Range("K5").Value = 1549.34052
$ws.Range("L5").Value = 3748.7142
$ws.Range("M5").Value = -1437.34052
$ws.Range("N5").Value = -3972.7142
$ws.Range("H122").Value = 1061
$ws.Range("I122").Value = 476.72726
$ws.Range("J122").Value = 1439.0588
$ws.Range("K122").Value = 4290.54534
$ws.Range("L122").Value = 12951.5292
$ws.Range("M122").Value = -1840.54534
$ws.Range("N122").Value = -17851.5292
$ws.Range("H123").Value = 3490.6206
$ws.Range("I123").Value = 1193
$ws.Range("J123").Value = 4699.8945
$ws.Range("K123").Value = 3579
$ws.Range("L123").Value = 14099.6835
$ws.Range("M123").Value = -1129
$ws.Range("N123").Value = -18999.6835
$ws.Range("H124").Value = 3742.9333
$ws.Range("I124").Value = 1065
$ws.Range("J124").Value = 4154.923
$ws.Range("K124").Value = 3195
$ws.Range("L124").Value = 12464.769
$ws.Range("M124").Value = 1715
$ws.Range("N124").Value = -22284.769
$ws.Range("H125").Value = 9457.048000000001
$ws.Range("I125").Value = 7000
$ws.Range("J125").Value = 9579.9
$ws.Range("K125").Value = 21000
$ws.Range("L125").Value = 28739.7
$ws.Range("M125").Value = -16080
$ws.Range("N125").Value = -38579.7
$ws.Range("H126").Value = 5199.8
$ws.Range("I126").Value = 4333
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 12999
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -8059
$ws.Range("N126").Value = -29380
$ws.Range("H134").Value = 1329.52
$ws.Range("I134").Value = 1369.2
$ws.Range("J134").Value = 1270
$ws.Range("K134").Value = 4107.6
$ws.Range("L134").Value = 3810
$ws.Range("M134").Value = 962.3999999999996
$ws.Range("N134").Value = -13950
$ws.Range("H135").Value = 684.7049
$ws.Range("I135").Value = 516.44684
$ws.Range("J135").Value = 1249.5714
$ws.Range("K135").Value = 4648.021559999999
$ws.Range("L135").Value = 11246.1426
$ws.Range("M135").Value = -2113.021559999999
$ws.Range("N135").Value = -16316.1426
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1299.3125
$ws.Range("I7").Value = 1183.7693
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 1183.7693
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -1071.7693
$ws.Range("N7").Value = -2024
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H126").Value = 1299.3125
$ws.Range("I126").Value = 1183.7693
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3551.3079
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -1081.3079
$ws.Range("N126").Value = -10340
$ws.Range("H136").Value = 7923
$ws.Range("I136").Value = 13605.2
$ws.Range("J136").Value = 3552.077
$ws.Range("K136").Value = 40815.60000000001
$ws.Range("L136").Value = 10656.231
$ws.Range("M136").Value = -38265.60000000001
$ws.Range("N136").Value = -15756.231
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 47830.68
$ws.Range("J123").Value = 47830.68
$ws.Range("L123").Value = 47830.68
$ws.Range("N123").Value = -57630.68
$ws.Range("H125").Value = 21485.334
$ws.Range("J125").Value = 21485.334
$ws.Range("L125").Value = 21485.334
$ws.Range("N125").Value = -31325.334
$ws.Range("H136").Value = 4124.5137
$ws.Range("I136").Value = 1039.5358
$ws.Range("J136").Value = 13722.223
$ws.Range("K136").Value = 3118.6074
$ws.Range("L136").Value = 41166.669
$ws.Range("M136").Value = -568.6074000000003
$ws.Range("N136").Value = -46266.669
